$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 27 de Julio de 2020 a las 21:19'

$ws.Range("B4").Value = 4407052
$ws.Range("C4").Value = 35213
$ws.Range("D4").Value = 2104800
$ws.Range("E4").Value = 2152122
$ws.Range("G4").Value = 282
$ws.Range("H4").Value = 150130

$ws.Range("B21").Value = 207071
$ws.Range("C21").Value = 330
$ws.Range("E21").Value = 7268

$ws.Range("D55").Value = 30900
$ws.Range("E55").Value = 1599

$ws.Range("A70").Value = 'Costa Rica'
$ws.Range("B70").Value = 15841
$ws.Range("C70").Value = 612
$ws.Range("D70").Value = 3824
$ws.Range("E70").Value = 11902
$ws.Range("G70").Value = 11
$ws.Range("H70").Value = 115

$ws.Range("A71").Value = 'Costa de Marfil'
$ws.Range("B71").Value = 15596
$ws.Range("D71").Value = 10178
$ws.Range("E71").Value = 5322
$ws.Range("H71").Value = 96

$ws.Range("A72").Value = 'Venezuela'
$ws.Range("B72").Value = 15463
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 9746
$ws.Range("E72").Value = 5575
$ws.Range("H72").Value = 142

$ws.Range("A73").Value = 'Chequia'
$ws.Range("B73").Value = 15421
$ws.Range("C73").Value = 97
$ws.Range("D73").Value = 11428
$ws.Range("E73").Value = 3620
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 373

$ws.Range("A118").Value = 'Cuba'
$ws.Range("B118").Value = 2532
$ws.Range("C118").Value = 37
$ws.Range("D118").Value = 2351
$ws.Range("E118").Value = 94
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 87

$ws.Range("A119").Value = 'Mali'
$ws.Range("B119").Value = 2513
$ws.Range("C119").Value = 3
$ws.Range("D119").Value = 1913
$ws.Range("E119").Value = 476
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 124

$ws.Range("A120").Value = 'Zimbabue'
$ws.Range("B120").Value = 2512
$ws.Range("D120").Value = 518
$ws.Range("E120").Value = 1960
$ws.Range("H120").Value = 34

$ws.Range("B121").Value = 2328
$ws.Range("C121").Value = 21
$ws.Range("D121").Value = 1550
$ws.Range("E121").Value = 756

$ws.Range("A122").Value = 'Suazilandia'
$ws.Range("B122").Value = 2316
$ws.Range("C122").Value = 109
$ws.Range("D122").Value = 1025
$ws.Range("E122").Value = 1257
$ws.Range("G122").Value = 2
$ws.Range("H122").Value = 34

$ws.Range("A123").Value = 'Sudan del Sur'
$ws.Range("B123").Value = 2305
$ws.Range("C123").Value = 43
$ws.Range("D123").Value = 1175
$ws.Range("E123").Value = 1084
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 46

$ws.Range("A134").Value = 'Mozambique'
$ws.Range("B134").Value = 1701
$ws.Range("C134").Value = 32
$ws.Range("D134").Value = 596
$ws.Range("E134").Value = 1094
$ws.Range("H134").Value = 11

$ws.Range("A135").Value = 'Yemen'
$ws.Range("B135").Value = 1691
$ws.Range("C135").Value = 10
$ws.Range("D135").Value = 833
$ws.Range("E135").Value = 375
$ws.Range("G135").Value = 4
$ws.Range("H135").Value = 483

$ws.Range("B145").Value = 1128
$ws.Range("C145").Value = 13
$ws.Range("D145").Value = 986
$ws.Range("E145").Value = 140

$ws.Range("A164").Value = 'Burundi'
$ws.Range("B164").Value = 378
$ws.Range("C164").Value = 17
$ws.Range("D164").Value = 301
$ws.Range("E164").Value = 76
$ws.Range("H164").Value = 1

$ws.Range("A165").Value = 'Guyana'
$ws.Range("B165").Value = 370
$ws.Range("D165").Value = 181
$ws.Range("E165").Value = 169
$ws.Range("H165").Value = 20

